# Update the "Förändrad" (Changed) date column (C) for rows 2-26
# from serial date 45233 (2023-11-03) to serial date 45243 (2023-11-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
